# GitHub Actions refresh of the cryptocurrency leaderboard: latest Price/Volume(1h)
# snapshot for every coin, plus three rows whose rank order changed (Bittensor now
# outranks EnergySwap, and ONDO dropped out of the top 50 in favor of Stacks).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the A1 cell reference and its new text. 'AsText' cells hold values
# (e.g. "1.00", "29.60") that Excel would otherwise auto-coerce to numbers, so the
# cell is briefly marked Text while the literal string is written, then restored to
# its original style so formatting stays untouched - matching the source feed, which
# stores every Price/Volume value as plain text.
$updates = @(
    @{ Cell = 'D2'; Value = '67.785.68'; AsText = $false }
    @{ Cell = 'E2'; Value = '  +0.65%  '; AsText = $false }
    @{ Cell = 'D3'; Value = '3.803.87'; AsText = $false }
    @{ Cell = 'E3'; Value = '  +0.90%  '; AsText = $false }
    @{ Cell = 'E4'; Value = '  +0.12%  '; AsText = $false }
    @{ Cell = 'D5'; Value = '596.76'; AsText = $true }
    @{ Cell = 'E5'; Value = '  +0.59%  '; AsText = $false }
    @{ Cell = 'D6'; Value = '167.25'; AsText = $true }
    @{ Cell = 'E6'; Value = '  +0.72%  '; AsText = $false }
    @{ Cell = 'E7'; Value = '  -0.04%  '; AsText = $false }
    @{ Cell = 'E8'; Value = '  +0.64%  '; AsText = $false }
    @{ Cell = 'D9'; Value = '0.161'; AsText = $true }
    @{ Cell = 'E9'; Value = '  +1.48%  '; AsText = $false }
    @{ Cell = 'E10'; Value = '  -1.10%  '; AsText = $false }
    @{ Cell = 'D11'; Value = '0.450'; AsText = $true }
    @{ Cell = 'E11'; Value = '  +0.35%  '; AsText = $false }
    @{ Cell = 'D12'; Value = '0.0000254'; AsText = $true }
    @{ Cell = 'E12'; Value = '  -0.35%  '; AsText = $false }
    @{ Cell = 'D13'; Value = '35.92'; AsText = $true }
    @{ Cell = 'E13'; Value = '  +0.32%  '; AsText = $false }
    @{ Cell = 'D14'; Value = '4.446.22'; AsText = $false }
    @{ Cell = 'E14'; Value = '  +1.04%  '; AsText = $false }
    @{ Cell = 'D15'; Value = '3.786.05'; AsText = $false }
    @{ Cell = 'E15'; Value = '  +0.21%  '; AsText = $false }
    @{ Cell = 'D16'; Value = '18.55'; AsText = $true }
    @{ Cell = 'E16'; Value = '  +4.31%  '; AsText = $false }
    @{ Cell = 'D17'; Value = '67.824.09'; AsText = $false }
    @{ Cell = 'E17'; Value = '  +0.74%  '; AsText = $false }
    @{ Cell = 'D18'; Value = '7.09'; AsText = $true }
    @{ Cell = 'E18'; Value = '  +1.86%  '; AsText = $false }
    @{ Cell = 'E19'; Value = '  +0.33%  '; AsText = $false }
    @{ Cell = 'D20'; Value = '461.81'; AsText = $true }
    @{ Cell = 'E20'; Value = '  +1.05%  '; AsText = $false }
    @{ Cell = 'D21'; Value = '9.91'; AsText = $true }
    @{ Cell = 'E21'; Value = '  -2.97%  '; AsText = $false }
    @{ Cell = 'D22'; Value = '0.702'; AsText = $true }
    @{ Cell = 'E22'; Value = '  +0.64%  '; AsText = $false }
    @{ Cell = 'E23'; Value = '  +0.97%  '; AsText = $false }
    @{ Cell = 'D24'; Value = '83.34'; AsText = $true }
    @{ Cell = 'E24'; Value = '  +0.13%  '; AsText = $false }
    @{ Cell = 'D25'; Value = '12.09'; AsText = $true }
    @{ Cell = 'E25'; Value = '  +2.21%  '; AsText = $false }
    @{ Cell = 'E26'; Value = '  -0.51%  '; AsText = $false }
    @{ Cell = 'E27'; Value = '  +0.06%  '; AsText = $false }
    @{ Cell = 'E28'; Value = '  +0.34%  '; AsText = $false }
    @{ Cell = 'D29'; Value = '3.947.95'; AsText = $false }
    @{ Cell = 'E29'; Value = '  +0.84%  '; AsText = $false }
    @{ Cell = 'E30'; Value = '  -0.08%  '; AsText = $false }
    @{ Cell = 'D31'; Value = '2.24'; AsText = $true }
    @{ Cell = 'E31'; Value = '  +2.93%  '; AsText = $false }
    @{ Cell = 'D32'; Value = '7.33'; AsText = $true }
    @{ Cell = 'E32'; Value = '  +2.05%  '; AsText = $false }
    @{ Cell = 'D33'; Value = '29.60'; AsText = $true }
    @{ Cell = 'D34'; Value = '0.999'; AsText = $true }
    @{ Cell = 'E34'; Value = '  -0.08%  '; AsText = $false }
    @{ Cell = 'D35'; Value = '9.07'; AsText = $true }
    @{ Cell = 'E35'; Value = '  -0.77%  '; AsText = $false }
    @{ Cell = 'D36'; Value = '3.745.76'; AsText = $false }
    @{ Cell = 'E36'; Value = '  +0.62%  '; AsText = $false }
    @{ Cell = 'D37'; Value = '0.100'; AsText = $true }
    @{ Cell = 'E37'; Value = '  +0.54%  '; AsText = $false }
    @{ Cell = 'E38'; Value = '  +2.48%  '; AsText = $false }
    @{ Cell = 'E39'; Value = '  +0.18%  '; AsText = $false }
    @{ Cell = 'D40'; Value = '1.00'; AsText = $true }
    @{ Cell = 'E40'; Value = '  +1.34%  '; AsText = $false }
    @{ Cell = 'D41'; Value = '5.78'; AsText = $true }
    @{ Cell = 'E41'; Value = '  +0.99%  '; AsText = $false }
    @{ Cell = 'D42'; Value = '1.00'; AsText = $true }
    @{ Cell = 'E42'; Value = '  +0.08%  '; AsText = $false }
    @{ Cell = 'E43'; Value = '  -0.01%  '; AsText = $false }
    @{ Cell = 'D44'; Value = '48.15'; AsText = $true }
    @{ Cell = 'E44'; Value = '  +2.78%  '; AsText = $false }
    @{ Cell = 'D45'; Value = '0.301'; AsText = $true }
    @{ Cell = 'E45'; Value = '  +1.48%  '; AsText = $false }
    @{ Cell = 'D46'; Value = '42.88'; AsText = $true }
    @{ Cell = 'E46'; Value = '  -1.61%  '; AsText = $false }
    @{ Cell = 'D47'; Value = '8.34'; AsText = $true }
    @{ Cell = 'E47'; Value = '  +0.02%  '; AsText = $false }
    @{ Cell = 'D48'; Value = '147.62'; AsText = $true }
    @{ Cell = 'E48'; Value = '  -0.15%  '; AsText = $false }
    @{ Cell = 'B49'; Value = 'Bittensor'; AsText = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; AsText = $false }
    @{ Cell = 'D49'; Value = '396.16'; AsText = $true }
    @{ Cell = 'E49'; Value = '  +1.37%  '; AsText = $false }
    @{ Cell = 'B50'; Value = 'EnergySwap'; AsText = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; AsText = $false }
    @{ Cell = 'D50'; Value = '27.07'; AsText = $true }
    @{ Cell = 'E50'; Value = '  +7.27%  '; AsText = $false }
    @{ Cell = 'B51'; Value = 'Stacks'; AsText = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; AsText = $false }
    @{ Cell = 'D51'; Value = '1.84'; AsText = $true }
    @{ Cell = 'E51'; Value = '  +1.26%  '; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.AsText) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates to $($ws.Name)"
